$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Construir un dispensador de alimento para mascotas con el fin de llenar el plato de comida con base las siguientes características; raza, edad, tamaño y peso de forma automática.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Construir un dispensador de alimento para mascotas con el fin de surtir de alimento al plato de comida con base las siguientes características; raza, edad, tamaño y peso de forma automática.",
    2
)
